$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" -> zero out a set of cells (row totals
# for several clients dropped to 0), and update the "X de 22"
# progress labels in row 24 to "0 de 22" for the matching columns.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("E3").Value = 0
$ws1.Range("G3").Value = 0
$ws1.Range("M3").Value = 0

$ws1.Range("M5").Value = 0

$ws1.Range("M6").Value = 0

$ws1.Range("D8").Value = 0
$ws1.Range("I8").Value = 0
$ws1.Range("M8").Value = 0

$ws1.Range("D10").Value = 0
$ws1.Range("L10").Value = 0

$ws1.Range("C11").Value = 0
$ws1.Range("E11").Value = 0
$ws1.Range("I11").Value = 0
$ws1.Range("M11").Value = 0
$ws1.Range("R11").Value = 0

$ws1.Range("M12").Value = 0

$ws1.Range("M13").Value = 0

$ws1.Range("M19").Value = 0

$ws1.Range("M23").Value = 0

$ws1.Range("C24").Value = "0 de 22"
$ws1.Range("D24").Value = "0 de 22"
$ws1.Range("E24").Value = "0 de 22"
$ws1.Range("G24").Value = "0 de 22"
$ws1.Range("I24").Value = "0 de 22"
$ws1.Range("L24").Value = "0 de 22"
$ws1.Range("M24").Value = "0 de 22"
$ws1.Range("R24").Value = "0 de 22"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL" -> roll the monthly columns forward one
# month (abril/mayo/junio/julio -> mayo/junio/julio/agosto), the
# oldest month's figures drop off and the new month starts at 0,
# plus a couple of figures get corrected along the way.
# Column F's width also narrows slightly (14 -> 12).
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Columns.Item(6).ColumnWidth = 11.17

$ws2.Range("C1").Value = "mayo"
$ws2.Range("D1").Value = "junio"
$ws2.Range("E1").Value = "julio"
$ws2.Range("F1").Value = "agosto"

$ws2.Range("C2").Value = -49.25
$ws2.Range("D2").Value = 0

$ws2.Range("C3").Value = 3054.27
$ws2.Range("D3").Value = 1317.8
$ws2.Range("E3").Value = 6836.54
$ws2.Range("F3").Value = 0

$ws2.Range("D4").Value = 777.8099999999999
$ws2.Range("E4").Value = 0

$ws2.Range("C5").Value = 226.8
$ws2.Range("D5").Value = 0
$ws2.Range("E5").Value = 366.83
$ws2.Range("F5").Value = 0

$ws2.Range("C6").Value = 2867.6
$ws2.Range("D6").Value = 25364.28
$ws2.Range("E6").Value = 2654.94
$ws2.Range("F6").Value = 0

$ws2.Range("D8").Value = 2261.64
$ws2.Range("E8").Value = 4423.73
$ws2.Range("F8").Value = 0

$ws2.Range("C9").Value = 2785.1
$ws2.Range("D9").Value = -22.68
$ws2.Range("E9").Value = 0

$ws2.Range("D10").Value = 851.4299999999999
$ws2.Range("E10").Value = 4229.93
$ws2.Range("F10").Value = 0

$ws2.Range("E11").Value = 5087.14
$ws2.Range("F11").Value = 0

$ws2.Range("C12").Value = 156.67
$ws2.Range("D12").Value = 11.52
$ws2.Range("E12").Value = 12246.22
$ws2.Range("F12").Value = 0

$ws2.Range("E13").Value = 7529.26
$ws2.Range("F13").Value = 0

$ws2.Range("E16").Value = 220.5

$ws2.Range("C17").Value = 0

$ws2.Range("C19").Value = 40.19
$ws2.Range("D19").Value = 0
$ws2.Range("E19").Value = 33.7
$ws2.Range("F19").Value = 0

$ws2.Range("C20").Value = 4277.03
$ws2.Range("D20").Value = 4321.33
$ws2.Range("E20").Value = 0

$ws2.Range("C23").Value = 24096.93
$ws2.Range("D23").Value = 4798.25
$ws2.Range("E23").Value = 16408.39
$ws2.Range("F23").Value = 0

$ws2.Range("C24").Value = 37455.34
$ws2.Range("D24").Value = 39681.38
$ws2.Range("E24").Value = 60037.18
$ws2.Range("F24").Value = 0
